{"js": "// FPLA-634: Minor c6 c6a word template updates\n//\n// 1. Remove the stray \"_GoBack\" bookmark that sat right after the\n//    \"<<courtName>>\" merge field.\n// 2. Re-word the \"special needs\" sentence:\n//      \"At the hearing, you can tell the court about any special needs\n//       or circumstances of the {child / children}.\"\n//    becomes\n//      \"At the hearing, you can tell the court if the child has any\n//       special needs or circumstances. \"\n//    and the \"_GoBack\" bookmark is re-inserted in the new sentence,\n//    right after \"...any special need\" (splitting the sentence into the\n//    same three runs the authored edit produced).\n\nconst body = context.document.body;\n\n// --- 1. Drop the old _GoBack bookmark (right after \"<<courtName>>\") ---\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// --- 2. Replace the sentence text --------------------------------------\nconst oldSentence =\n  \"At the hearing, you can tell the court about any special needs or circumstances of the {child / children}.\";\nconst newSentence =\n  \"At the hearing, you can tell the court if the child has any special needs or circumstances. \";\n\nconst sentenceHits = body.search(oldSentence, { matchCase: true });\nsentenceHits.load(\"items\");\nawait context.sync();\n\nif (sentenceHits.items.length === 0) {\n  throw new Error(\"Could not find the target sentence to update.\");\n}\n\nsentenceHits.items[0].insertText(newSentence, \"Replace\");\nawait context.sync();\n\n// --- 3. Re-create the run split the original edit left behind ----------\n// First split off \"At the hearing, you can tell the court \" from the rest\n// using a throw-away bookmark purely as a run-boundary marker.\nconst firstPartHits = body.search(\n  \"At the hearing, you can tell the court \",\n  { matchCase: true }\n);\nfirstPartHits.load(\"items\");\nawait context.sync();\n\nif (firstPartHits.items.length === 0) {\n  throw new Error(\"Could not find the first part of the new sentence.\");\n}\n\nfirstPartHits.items[0].getRange(\"After\").insertBookmark(\"__tmp_split__\");\nawait context.sync();\n\n// Then place the real \"_GoBack\" bookmark right after \"...any special\n// need\" (before the trailing \"s\" of \"needs\"), which splits off the\n// second and third runs.\nconst secondPartHits = body.search(\n  \"if the child has any special need\",\n  { matchCase: true }\n);\nsecondPartHits.load(\"items\");\nawait context.sync();\n\nif (secondPartHits.items.length === 0) {\n  throw new Error(\"Could not find the second part of the new sentence.\");\n}\n\nsecondPartHits.items[0].getRange(\"After\").insertBookmark(\"_GoBack\");\nawait context.sync();\n\n// Remove the temporary marker bookmark now that the runs are split; the\n// split survives the bookmark's removal.\ncontext.document.deleteBookmark(\"__tmp_split__\");\nawait context.sync();\n", "ps1": "# FPLA-634: Minor c6 c6a word template updates\n#\n# 1. Remove the stray \"_GoBack\" bookmark that sat right after the\n#    \"<<courtName>>\" merge field.\n# 2. Re-word the \"special needs\" sentence:\n#      \"At the hearing, you can tell the court about any special needs\n#       or circumstances of the {child / children}.\"\n#    becomes\n#      \"At the hearing, you can tell the court if the child has any\n#       special needs or circumstances. \"\n#    and the \"_GoBack\" bookmark is re-inserted in the new sentence,\n#    right after \"...any special need\" (splitting the sentence into the\n#    same three runs the authored edit produced).\n\n$d = $word.ActiveDocument\n\n# --- 1. Drop the old _GoBack bookmark (right after \"<<courtName>>\") ---\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# --- 2. Replace the sentence text --------------------------------------\n$oldSentence = \"At the hearing, you can tell the court about any special needs or circumstances of the {child / children}.\"\n$newSentence = \"At the hearing, you can tell the court if the child has any special needs or circumstances. \"\n\n$sentenceRange = $d.Content\n$sentenceFind = $sentenceRange.Find\n$sentenceFind.ClearFormatting()\n$found = $sentenceFind.Execute($oldSentence)\nif (-not $found) {\n    throw \"Could not find the target sentence to update.\"\n}\n$sentenceRange.Text = $newSentence\n\n# --- 3. Re-create the run split the original edit left behind ----------\n# First split off \"At the hearing, you can tell the court \" from the rest\n# using a throw-away bookmark purely as a run-boundary marker.\n$firstPartRange = $d.Content\n$firstPartFind = $firstPartRange.Find\n$firstPartFind.ClearFormatting()\n$foundFirst = $firstPartFind.Execute(\"At the hearing, you can tell the court \")\nif (-not $foundFirst) {\n    throw \"Could not find the first part of the new sentence.\"\n}\n$firstPartRange.Collapse(0)  # wdCollapseEnd\n$d.Bookmarks.Add(\"__tmp_split__\", $firstPartRange) | Out-Null\n\n# Then place the real \"_GoBack\" bookmark right after \"...any special\n# need\" (before the trailing \"s\" of \"needs\"), which splits off the\n# second and third runs.\n$secondPartRange = $d.Content\n$secondPartFind = $secondPartRange.Find\n$secondPartFind.ClearFormatting()\n$foundSecond = $secondPartFind.Execute(\"if the child has any special need\")\nif (-not $foundSecond) {\n    throw \"Could not find the second part of the new sentence.\"\n}\n$secondPartRange.Collapse(0)  # wdCollapseEnd\n$d.Bookmarks.Add(\"_GoBack\", $secondPartRange) | Out-Null\n\n# Remove the temporary marker bookmark now that the runs are split; the\n# split survives the bookmark's removal.\n$d.Bookmarks(\"__tmp_split__\").Delete()\n"}
